$d = $word.ActiveDocument
$t = $d.Tables(1)

$newValues = @(
    "28+59=87",
    "10+11=21",
    "81-64=17",
    "55-21=34",
    "88-82=6",
    "26+50=76",
    "69+28=97",
    "73+7=80",
    "90-27=63",
    "74+8=82",
    "31-3=28",
    "63+17=80",
    "60-8=52",
    "38+14=52",
    "22+47=69",
    "27+26=53",
    "28-16=12",
    "16+48=64",
    "43+44=87",
    "61-41=20",
    "13+28=41",
    "13+43=56",
    "89-26=63",
    "53+43=96",
    "48+47=95",
    "4+30=34",
    "8-6=2",
    "47-35=12",
    "36+23=59",
    "82-39=43",
    "46-16=30",
    "69+28=97",
    "47+39=86",
    "70-42=28",
    "16+43=59",
    "8+89=97",
    "72+18=90",
    "43+23=66",
    "82-26=56",
    "19-13=6",
    "68+2=70",
    "64+21=85",
    "68-31=37",
    "60-59=1",
    "1+49=50",
    "14+7=21",
    "17+58=75",
    "79-77=2",
    "42+32=74",
    "30+1=31",
    "97-48=49",
    "47+23=70",
    "16+68=84",
    "4+24=28",
    "43+41=84",
    "25+15=40",
    "47+25=72",
    "73-63=10",
    "65-43=22",
    "86+8=94",
    "42+10=52",
    "72-5=67",
    "33-18=15",
    "43-10=33",
    "64-45=19",
    "9-8=1",
    "32+43=75",
    "25-10=15",
    "55+7=62",
    "95-71=24",
    "43-27=16",
    "22+9=31",
    "67-11=56",
    "36+16=52",
    "4+94=98",
    "19-1=18",
    "10+69=79",
    "91-90=1",
    "41-0=41",
    "75-65=10",
    "83+0=83",
    "75-2=73",
    "52+35=87",
    "80-52=28",
    "35-6=29",
    "46+52=98",
    "74-58=16",
    "35+23=58",
    "42+2=44",
    "80-14=66",
    "13+19=32",
    "65-26=39",
    "36+59=95",
    "84-12=72",
    "77+8=85",
    "17+73=90",
    "73-23=50",
    "61-58=3",
    "65-20=45",
    "24+7=31"
)

$rows = 20
$cols = 5
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $rng = $cell.Range
        $rng.End = $rng.End - 1
        $rng.Text = $newValues[$idx]
        $idx++
    }
}

Write-Host "Updated" $idx "cells"
